# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 43 (pushing the existing rows 43-83
# down to 44-84) and populate it with the new "Ají" observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 43; all rows below shift down by one.
$ws.Rows("43:43").Insert()

# Populate the newly inserted row 43 with the new record.
$ws.Range("A43").Value = 7
$ws.Range("B43").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C43").Value = "Ñuble"
$ws.Range("D43").Value = 44673
$ws.Range("E43").Value = 16
$ws.Range("F43").Value = 100112021
$ws.Range("G43").Value = "Ají"
$ws.Range("H43").Value = "Cacho cabra verde"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 30
$ws.Range("K43").Value = 16000
$ws.Range("L43").Value = 17000
$ws.Range("M43").Value = 16500
$ws.Range("N43").Value = "$/caja 25 kilos"
$ws.Range("O43").Value = "Provincia de Diguillín"
$ws.Range("P43").Value = 660
$ws.Range("Q43").Value = 25
$ws.Range("R43").Value = "Hortaliza"

# Match the date formatting/style used by the rest of the column D cells.
$ws.Range("D43").NumberFormat = $ws.Range("D44").NumberFormat()
